# hero_board_bot_annotated.pptx - "Highlight which parts must be hand assembled"
#
# 1) Refresh the cached datetimeFigureOut text (master + all 11 layouts):
#    2020/01/26 -> 2020-05-20
# 2) Nudge the C/D/E label boxes (and their leader-line connectors) to sit
#    better over the board photo.
# 3) Give the G / H / N label boxes a translucent yellow fill so the
#    hand-assembled connectors stand out.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder text on the slide master and every slide layout.
# ---------------------------------------------------------------------
$master = $p.Slides.Item(1).Master
$master.Shapes.Item("Date Placeholder 3").TextFrame.TextRange.Text = "2020-05-20"

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    $layout.Shapes.Item("Date Placeholder 3").TextFrame.TextRange.Text = "2020-05-20"
}

# ---------------------------------------------------------------------
# 2) Reposition C / D / E boxes and their connector arrows.
#    (Left/Top are single-precision points under the hood, so the
#    literal values below are pre-nudged to land exactly on the target
#    EMU offset after the point->EMU round-trip.)
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

$rectC = $s.Shapes.Item("Rectangle 37")
$rectC.Left = 504.62451171875
$rectC.Top = 28.26401710510254

$connC = $s.Shapes.Item("Straight Arrow Connector 38")
$connC.Left = 530.1115112304688
$connC.Top = 75.28441619873047

$rectD = $s.Shapes.Item("Rectangle 40")
$rectD.Left = 575.12939453125
$rectD.Top = 28.21685218811035

$rectE = $s.Shapes.Item("Rectangle 41")
$rectE.Left = 644.3032836914062
$rectE.Top = 28.21685218811035

$connFromD = $s.Shapes.Item("Straight Arrow Connector 47")
$connFromD.Left = 590.01123046875
$connFromD.Top = 75.23725128173828

$connFromE = $s.Shapes.Item("Straight Arrow Connector 49")
$connFromE.Left = 649.8165893554688
$connFromE.Top = 75.23725128173828

# ---------------------------------------------------------------------
# 3) Highlight the hand-assembled connector labels: G, H, N.
# ---------------------------------------------------------------------
$highlightNames = @("Rectangle 46", "Rectangle 64", "Rectangle 81")
foreach ($name in $highlightNames) {
    $sh = $s.Shapes.Item($name)
    $sh.Fill.Visible = $true
    $sh.Fill.Solid()
    $sh.Fill.ForeColor.RGB = 65535
    $sh.Fill.Transparency = 0.75
}
